# Static_Offset.xlsx update — "Updating to the latest CPRA slide deck workflow."
#
# 1. Row 2 (A2/B2) changes from a USACE/numeric-gage-id row to a USGS row
#    with a text gage id ("073802332").
# 2. Row 7's B cell (gage id) switches from the numeric 82742 to the text
#    gage id "82740".
# 3. A new data row is inserted at row 23 (USACE / 76305 / 0.5), pushing the
#    former rows 23-29 down to 24-30.
# 4. Several offset values in column C are updated from 0 to their new
#    values (including some now-negative offsets).
# 5. Selection/active cell moves to C24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: USACE -> USGS, numeric gage id -> text gage id "073802332" ---
$ws.Range("A2").Value = "USGS"
$ws.Range("B2").Value = "073802332"

# --- Row 7: numeric gage id -> text gage id "82740" ---
$ws.Range("B7").Value = "82740"

# --- Insert a brand new row at position 23 (shifts old rows 23-29 -> 24-30) ---
$ws.Rows(23).Insert()

# New row 23 content: A=USACE (shared-string), B=76305 as a genuine number
# (not text) even though the column's style uses a Text number format, C=0.5
$ws.Range("A23").Value = "USACE"
$ws.Range("B23").NumberFormat = "General"
$ws.Range("B23").Value = 76305
$ws.Range("B23").NumberFormat = "@"
$ws.Range("C23").Value = 0.5

# --- Column C offset updates (row numbers below are POST-insert numbering) ---
$ws.Range("C1").Value = 0.25
$ws.Range("C2").Value = 0.25
$ws.Range("C4").Value = 0.75
$ws.Range("C9").Value = 0.25
$ws.Range("C11").Value = 0.75
$ws.Range("C13").Value = 0.3
$ws.Range("C15").Value = 0.25
$ws.Range("C16").Value = 0.3
$ws.Range("C17").Value = 0.5
$ws.Range("C19").Value = 0.75
$ws.Range("C20").Value = 0.6
$ws.Range("C21").Value = 0.5
$ws.Range("C22").Value = 1.5
$ws.Range("C24").Value = -1.25
$ws.Range("C25").Value = -1.25
$ws.Range("C26").Value = -0.5
$ws.Range("C28").Value = -0.25
$ws.Range("C30").Value = -1.5

# --- Selection moves to C24 ---
$ws.Range("C24").Select()
